$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) columns store plain-decimal-looking text
# (e.g. "69.10", "213.29", "29.831.02"). Force the whole data range to Text
# format before writing so Excel does not auto-convert values that look like
# numbers (which would silently drop significant trailing zeros / dots).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.839.04'
$ws.Range('D3').Value = '1.619.04'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  -0.74%  '
$ws.Range('D5').Value = '213.29'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('D8').Value = '29.11'
$ws.Range('E8').Value = '  +8.92%  '
$ws.Range('E9').Value = '  +3.17%  '
$ws.Range('D10').Value = '0.0606'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.849.43'
$ws.Range('D13').Value = '1.612.69'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').Value = '0.566'
$ws.Range('E14').Value = '  +5.85%  '
$ws.Range('D15').Value = '3.90'
$ws.Range('E15').Value = '  +5.11%  '
$ws.Range('D16').Value = '29.852.47'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '8.85'
$ws.Range('E17').Value = '  +16.02%  '
$ws.Range('D18').Value = '64.33'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').Value = '241.34'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('D23').Value = '9.58'
$ws.Range('E23').Value = '  +4.01%  '
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').Value = '155.20'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('D26').Value = '15.60'
$ws.Range('E26').Value = '  +2.24%  '
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  +2.93%  '
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('E31').Value = '  +5.25%  '
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  +3.69%  '
$ws.Range('D33').Value = '3.21'
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('D34').Value = '1.414.97'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '1.64'
$ws.Range('E35').Value = '  +6.40%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '2.89'
$ws.Range('E36').Value = '  +2.44%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.03'
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').Value = '0.0169'
$ws.Range('D40').Value = '0.554'
$ws.Range('E40').Value = '  +3.12%  '
$ws.Range('D41').Value = '0.0501'
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('D42').Value = '0.826'
$ws.Range('E42').Value = '  +3.63%  '
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').Value = '69.10'
$ws.Range('E44').Value = '  +4.96%  '
$ws.Range('D45').Value = '53.34'
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('E46').Value = '  +19.00%  '
$ws.Range('D47').Value = '0.993'
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('E48').Value = '  +2.93%  '
$ws.Range('D49').Value = '1.759.38'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('D50').Value = '88.35'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').Value = '0.0₆0110'
$ws.Range('E51').Value = '  +5.11%  '
